$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts D:K -> E:L), carrying values/styles right
$ws.Columns("D:D").Insert()

# Copy number formats/styles from column E into the newly inserted column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# Rows 5, 6, 37 and 79 are section headers with a single label cell and never
# had a column D entry - undo the blank styled cell the paste above created.
$ws.Range("D5").Clear()
$ws.Range("D6").Clear()
$ws.Range("D37").Clear()
$ws.Range("D79").Clear()

# Match the new column's width to its neighbours (D:K are a uniform block)
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# Populate new column D with the newest period (FY2018) figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1207800
$ws.Range("D9").Value = 956700
$ws.Range("D10").Value = 251100
$ws.Range("D12").Value = 15200
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 5600
$ws.Range("D15").Value = 2300
$ws.Range("D17").Value = 1146300
$ws.Range("D18").Value = 61500
$ws.Range("D20").Value = -42700
$ws.Range("D21").Value = 54300
$ws.Range("D22").Value = 2500
$ws.Range("D23").Value = 16300
$ws.Range("D24").Value = 6600
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 9700
$ws.Range("D27").Value = 9700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 11100
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 42700
$ws.Range("D33").Value = 20800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 20800
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 70600
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 130500
$ws.Range("D44").Value = 214900
$ws.Range("D45").Value = 23300
$ws.Range("D46").Value = 439400
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 251000
$ws.Range("D49").Value = 97100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 12900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 800300
$ws.Range("D57").Value = 49600
$ws.Range("D58").Value = 800
$ws.Range("D59").Value = 89300
$ws.Range("D60").Value = 139800
$ws.Range("D61").Value = 17300
$ws.Range("D62").Value = 89400
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 246400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 548400
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 553900
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 20800
$ws.Range("D83").Value = 35500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 76400
$ws.Range("D91").Value = -27700
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -33800
$ws.Range("D96").Value = -8400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -13600
$ws.Range("D101").Value = -100
$ws.Range("D102").Value = 28800

# A handful of prior-year comparatives were restated in this edit
$ws.Range("E9").Value = 926600
$ws.Range("F9").Value = 784700
$ws.Range("E10").Value = 212800
$ws.Range("F10").Value = 184500
$ws.Range("E14").Value = 900
$ws.Range("E17").Value = 1099400
$ws.Range("F17").Value = 940400
$ws.Range("E18").Value = 40000
$ws.Range("F18").Value = 28900
$ws.Range("E20").Value = -1500
$ws.Range("F20").Value = -1800
$ws.Range("E32").Value = 1500
$ws.Range("F32").Value = 1800
